$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.333.93'
$ws.Range("E2").Value = '  -4.67%  '
$ws.Range("D3").Value = '3.263.95'
$ws.Range("E3").Value = '  -7.11%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("D5").Value = '''593.54'
$ws.Range("E5").Value = '  -4.98%  '
$ws.Range("D6").Value = '''151.17'
$ws.Range("E6").Value = '  -12.12%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.255.71'
$ws.Range("E8").Value = '  -7.20%  '
$ws.Range("D9").Value = '''0.542'
$ws.Range("E9").Value = '  -10.98%  '
$ws.Range("E10").Value = '  -14.36%  '
$ws.Range("E11").Value = '  -7.64%  '
$ws.Range("D12").Value = '''0.511'
$ws.Range("E12").Value = '  -12.42%  '
$ws.Range("D13").Value = '''38.14'
$ws.Range("E13").Value = '  -17.45%  '
$ws.Range("E14").Value = '  -11.58%  '
$ws.Range("D15").Value = '3.787.39'
$ws.Range("E15").Value = '  -7.23%  '
$ws.Range("D16").Value = '67.328.45'
$ws.Range("E16").Value = '  -4.84%  '
$ws.Range("D17").Value = '3.266.72'
$ws.Range("E17").Value = '  -7.10%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").Value = '''0.114'
$ws.Range("E18").Value = '  -6.34%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = '''7.22'
$ws.Range("E19").Value = '  -14.21%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '''533.06'
$ws.Range("E20").Value = '  -11.95%  '
$ws.Range("D21").Value = '''15.07'
$ws.Range("E21").Value = '  -14.87%  '
$ws.Range("E22").Value = '  -13.33%  '
$ws.Range("D23").Value = '''7.90'
$ws.Range("E23").Value = '  -13.05%  '
$ws.Range("D24").Value = '''85.48'
$ws.Range("E24").Value = '  -11.94%  '
$ws.Range("D25").Value = '''13.58'
$ws.Range("E25").Value = '  -12.27%  '
$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  -12.84%  '
$ws.Range("E28").Value = '  -12.24%  '
$ws.Range("D29").Value = '''8.02'
$ws.Range("E29").Value = '  -10.87%  '
$ws.Range("E30").Value = '  -16.62%  '
$ws.Range("E31").Value = '  -11.51%  '
$ws.Range("E32").Value = '  -11.09%  '
$ws.Range("D33").Value = '''543.09'
$ws.Range("E33").Value = '  -12.66%  '
$ws.Range("E34").Value = '  -17.83%  '
$ws.Range("D35").Value = '''5.71'
$ws.Range("E35").Value = '  -15.82%  '
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("E37").Value = '  -7.48%  '
$ws.Range("D38").Value = '''53.16'
$ws.Range("E38").Value = '  -6.16%  '
$ws.Range("D39").Value = '''0.0858'
$ws.Range("E39").Value = '  -13.61%  '
$ws.Range("E40").Value = '  -10.03%  '
$ws.Range("D41").Value = '''9.08'
$ws.Range("E41").Value = '  -15.99%  '
$ws.Range("D42").Value = '''2.71'
$ws.Range("D43").Value = '2.928.58'
$ws.Range("E43").Value = '  -12.18%  '
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").Value = '''0.262'
$ws.Range("E44").Value = '  -15.52%  '
$ws.Range("B45").Value = 'PEPE'
$ws.Range("C45").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D45").Value = '0.0₃0588'
$ws.Range("E45").Value = '  -18.41%  '
$ws.Range("D46").Value = '''26.94'
$ws.Range("E46").Value = '  -15.22%  '
$ws.Range("E47").Value = '  -13.90%  '
$ws.Range("D49").Value = '''126.92'
$ws.Range("E49").Value = '  -5.24%  '
$ws.Range("D50").Value = '''2.33'
$ws.Range("E50").Value = '  -20.47%  '
$ws.Range("E51").Value = '  -12.81%  '
